$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fn1"
$ws.Range("C2").Value = "Itga4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 27.03890566666666
$ws.Range("H2").Value = 81.116717
$ws.Range("I2").Value = 0.07096188219033728
$ws.Range("J2").Value = 0.07096188219033729
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.04517333333333
$ws.Range("N2").Value = 78.13552
$ws.Range("O2").Value = 0.9210237118384171
$ws.Range("P2").Value = 0.921023711838417
$ws.Range("Q2").Value = 704.2329848319822
$ws.Range("R2").Value = 6338.096863487839
$ws.Range("S2").Value = 0.06535757613398491
$ws.Range("T2").Value = 0.06535757613398491

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fn1"
$ws.Range("C3").Value = "Itga4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 27.03890566666666
$ws.Range("H3").Value = 81.116717
$ws.Range("I3").Value = 0.07096188219033728
$ws.Range("J3").Value = 0.07096188219033729
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.3302223333333333
$ws.Range("N3").Value = 0.990667
$ws.Range("O3").Value = 0.01167750336256582
$ws.Range("P3").Value = 0.01167750336256582
$ws.Range("Q3").Value = 8.928850520026556
$ws.Range("R3").Value = 80.35965468023899
$ws.Range("S3").Value = 0.0008286576178916634
$ws.Range("T3").Value = 0.0008286576178916635

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fn1"
$ws.Range("C4").Value = "Itga4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 27.03890566666666
$ws.Range("H4").Value = 81.116717
$ws.Range("I4").Value = 0.07096188219033728
$ws.Range("J4").Value = 0.07096188219033729
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.903109
$ws.Range("N4").Value = 5.709327
$ws.Range("O4").Value = 0.06729878479901708
$ws.Range("P4").Value = 0.06729878479901708
$ws.Range("Q4").Value = 51.45798472438432
$ws.Range("R4").Value = 463.1218625194589
$ws.Range("S4").Value = 0.004775648438460711
$ws.Range("T4").Value = 0.004775648438460712

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fn1"
$ws.Range("C5").Value = "Itga4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 345.566579
$ws.Range("H5").Value = 1036.699737
$ws.Range("I5").Value = 0.9069174311350353
$ws.Range("J5").Value = 0.9069174311350354
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 26.04517333333333
$ws.Range("N5").Value = 78.13552
$ws.Range("O5").Value = 0.9210237118384171
$ws.Range("P5").Value = 0.921023711838417
$ws.Range("Q5").Value = 9000.341448262026
$ws.Range("R5").Value = 81003.07303435824
$ws.Range("S5").Value = 0.8352924587549523
$ws.Range("T5").Value = 0.8352924587549523

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fn1"
$ws.Range("C6").Value = "Itga4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 345.566579
$ws.Range("H6").Value = 1036.699737
$ws.Range("I6").Value = 0.9069174311350353
$ws.Range("J6").Value = 0.9069174311350354
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3302223333333333
$ws.Range("N6").Value = 0.990667
$ws.Range("O6").Value = 0.01167750336256582
$ws.Range("P6").Value = 0.01167750336256582
$ws.Range("Q6").Value = 114.1138020393977
$ws.Range("R6").Value = 1027.024218354579
$ws.Range("S6").Value = 0.01059053135164893
$ws.Range("T6").Value = 0.01059053135164893

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fn1"
$ws.Range("C7").Value = "Itga4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 345.566579
$ws.Range("H7").Value = 1036.699737
$ws.Range("I7").Value = 0.9069174311350353
$ws.Range("J7").Value = 0.9069174311350354
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.903109
$ws.Range("N7").Value = 5.709327
$ws.Range("O7").Value = 0.06729878479901708
$ws.Range("P7").Value = 0.06729878479901708
$ws.Range("Q7").Value = 657.650866594111
$ws.Range("R7").Value = 5918.857799346999
$ws.Range("S7").Value = 0.06103444102843413
$ws.Range("T7").Value = 0.06103444102843414

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fn1"
$ws.Range("C8").Value = "Itga4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 8.428738666666668
$ws.Range("H8").Value = 25.286216
$ws.Range("I8").Value = 0.0221206866746274
$ws.Range("J8").Value = 0.02212068667462741
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 26.04517333333333
$ws.Range("N8").Value = 78.13552
$ws.Range("O8").Value = 0.9210237118384171
$ws.Range("P8").Value = 0.921023711838417
$ws.Range("Q8").Value = 219.5279595547022
$ws.Range("R8").Value = 1975.75163599232
$ws.Range("S8").Value = 0.02037367694947994
$ws.Range("T8").Value = 0.02037367694947994

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fn1"
$ws.Range("C9").Value = "Itga4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 8.428738666666668
$ws.Range("H9").Value = 25.286216
$ws.Range("I9").Value = 0.0221206866746274
$ws.Range("J9").Value = 0.02212068667462741
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.3302223333333333
$ws.Range("N9").Value = 0.990667
$ws.Range("O9").Value = 0.01167750336256582
$ws.Range("P9").Value = 0.01167750336256582
$ws.Range("Q9").Value = 2.783357749563556
$ws.Range("R9").Value = 25.050219746072
$ws.Range("S9").Value = 0.0002583143930252265
$ws.Range("T9").Value = 0.0002583143930252266

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fn1"
$ws.Range("C10").Value = "Itga4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 8.428738666666668
$ws.Range("H10").Value = 25.286216
$ws.Range("I10").Value = 0.0221206866746274
$ws.Range("J10").Value = 0.02212068667462741
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.903109
$ws.Range("N10").Value = 5.709327
$ws.Range("O10").Value = 0.06729878479901708
$ws.Range("P10").Value = 0.06729878479901708
$ws.Range("Q10").Value = 16.04080841518133
$ws.Range("R10").Value = 144.367275736632
$ws.Range("S10").Value = 0.001488695332122234
$ws.Range("T10").Value = 0.001488695332122234
